$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.366139666666667
$ws.Range("H2").Value = 10.098419
$ws.Range("I2").Value = 0.01725116351498256
$ws.Range("J2").Value = 0.01815407111703398
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.546648333333333
$ws.Range("N2").Value = 7.639944999999999
$ws.Range("O2").Value = 0.01458319278139283
$ws.Range("P2").Value = 0.01462456071422719
$ws.Range("Q2").Value = 8.572373971883888
$ws.Range("R2").Value = 77.15136574695499
$ws.Range("S2").Value = 0.000251577043242321
$ws.Range("T2").Value = 0.0002654953152614616
$ws.Range("G3").Value = 3.366139666666667
$ws.Range("H3").Value = 10.098419
$ws.Range("I3").Value = 0.01725116351498256
$ws.Range("J3").Value = 0.01815407111703398
$ws.Range("O3").Value = 0.01449435301115475
$ws.Range("P3").Value = 0.01453546893349296
$ws.Range("Q3").Value = 8.520151681095223
$ws.Range("R3").Value = 76.681365129857
$ws.Range("S3").Value = 0.0002500444538393104
$ws.Range("T3").Value = 0.0002638779367380693
$ws.Range("G4").Value = 3.366139666666667
$ws.Range("H4").Value = 10.098419
$ws.Range("I4").Value = 0.01725116351498256
$ws.Range("J4").Value = 0.01815407111703398
$ws.Range("M4").Value = 94.63104
$ws.Range("N4").Value = 283.89312
$ws.Range("O4").Value = 0.5418976312357076
$ws.Range("P4").Value = 0.5434348244380536
$ws.Range("Q4").Value = 318.54129744192
$ws.Range("R4").Value = 2866.87167697728
$ws.Range("S4").Value = 0.009348364644828912
$ws.Range("T4").Value = 0.009865554450321298
$ws.Range("G5").Value = 3.366139666666667
$ws.Range("H5").Value = 10.098419
$ws.Range("I5").Value = 0.01725116351498256
$ws.Range("J5").Value = 0.01815407111703398
$ws.Range("M5").Value = 1.481899
$ws.Range("N5").Value = 2.963798
$ws.Range("O5").Value = 0.00848598470259403
$ws.Range("P5").Value = 0.005673371182083786
$ws.Range("Q5").Value = 4.988279005893667
$ws.Range("R5").Value = 29.929674035362
$ws.Range("S5").Value = 0.0001463931096900903
$ws.Range("T5").Value = 0.0001029947839128802
$ws.Range("G6").Value = 3.366139666666667
$ws.Range("H6").Value = 10.098419
$ws.Range("I6").Value = 0.01725116351498256
$ws.Range("J6").Value = 0.01815407111703398
$ws.Range("M6").Value = 73.43827566666666
$ws.Range("N6").Value = 220.314827
$ws.Range("O6").Value = 0.420538838269151
$ws.Range("P6").Value = 0.4217317747321426
$ws.Range("Q6").Value = 247.2034927731681
$ws.Range("R6").Value = 2224.831434958513
$ws.Range("S6").Value = 0.007254784263381929
$ws.Range("T6").Value = 0.007656148630800269
$ws.Range("I7").Value = 0.7504462978934635
$ws.Range("J7").Value = 0.7897238612132288
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.546648333333333
$ws.Range("N7").Value = 7.639944999999999
$ws.Range("O7").Value = 0.01458319278139283
$ws.Range("P7").Value = 0.01462456071422719
$ws.Range("Q7").Value = 372.9085464740639
$ws.Range("R7").Value = 3356.176918266575
$ws.Range("S7").Value = 0.01094390303426293
$ws.Range("T7").Value = 0.01154936455578679
$ws.Range("I8").Value = 0.7504462978934635
$ws.Range("J8").Value = 0.7897238612132288
$ws.Range("O8").Value = 0.01449435301115475
$ws.Range("P8").Value = 0.01453546893349296
$ws.Range("S8").Value = 0.01087723355758205
$ws.Range("T8").Value = 0.01147900665070299
$ws.Range("I9").Value = 0.7504462978934635
$ws.Range("J9").Value = 0.7897238612132288
$ws.Range("M9").Value = 94.63104
$ws.Range("N9").Value = 283.89312
$ws.Range("O9").Value = 0.5418976312357076
$ws.Range("P9").Value = 0.5434348244380536
$ws.Range("Q9").Value = 13856.9283853728
$ws.Range("R9").Value = 124712.3554683552
$ws.Range("S9").Value = 0.406665071198074
$ws.Range("T9").Value = 0.4291634478729528
$ws.Range("I10").Value = 0.7504462978934635
$ws.Range("J10").Value = 0.7897238612132288
$ws.Range("M10").Value = 1.481899
$ws.Range("N10").Value = 2.963798
$ws.Range("O10").Value = 0.00848598470259403
$ws.Range("P10").Value = 0.005673371182083786
$ws.Range("Q10").Value = 216.9961179477217
$ws.Range("R10").Value = 1301.97670768633
$ws.Range("S10").Value = 0.006368275804042253
$ws.Range("T10").Value = 0.004480396596011068
$ws.Range("I11").Value = 0.7504462978934635
$ws.Range("J11").Value = 0.7897238612132288
$ws.Range("M11").Value = 73.43827566666666
$ws.Range("N11").Value = 220.314827
$ws.Range("O11").Value = 0.420538838269151
$ws.Range("P11").Value = 0.4217317747321426
$ws.Range("Q11").Value = 10753.64834475312
$ws.Range("R11").Value = 96782.83510277804
$ws.Range("S11").Value = 0.3155918142995024
$ws.Range("T11").Value = 0.3330516455377752
$ws.Range("G12").Value = 14.89209833333333
$ws.Range("H12").Value = 44.676295
$ws.Range("I12").Value = 0.07632066665966204
$ws.Range("J12").Value = 0.08031520940808551
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.546648333333333
$ws.Range("N12").Value = 7.639944999999999
$ws.Range("O12").Value = 0.01458319278139283
$ws.Range("P12").Value = 0.01462456071422719
$ws.Range("Q12").Value = 37.92493740041944
$ws.Range("R12").Value = 341.3244366037749
$ws.Range("S12").Value = 0.001112998995102272
$ws.Range("T12").Value = 0.001174574656264417
$ws.Range("G13").Value = 14.89209833333333
$ws.Range("H13").Value = 44.676295
$ws.Range("I13").Value = 0.07632066665966204
$ws.Range("J13").Value = 0.08031520940808551
$ws.Range("O13").Value = 0.01449435301115475
$ws.Range("P13").Value = 0.01453546893349296
$ws.Range("Q13").Value = 37.69390138687611
$ws.Range("R13").Value = 339.245112481885
$ws.Range("S13").Value = 0.00110621868461181
$ws.Range("T13").Value = 0.001167419231238209
$ws.Range("G14").Value = 14.89209833333333
$ws.Range("H14").Value = 44.676295
$ws.Range("I14").Value = 0.07632066665966204
$ws.Range("J14").Value = 0.08031520940808551
$ws.Range("M14").Value = 94.63104
$ws.Range("N14").Value = 283.89312
$ws.Range("O14").Value = 0.5418976312357076
$ws.Range("P14").Value = 0.5434348244380536
$ws.Range("Q14").Value = 1409.2547530656
$ws.Range("R14").Value = 12683.2927775904
$ws.Range("S14").Value = 0.0413579884772009
$ws.Range("T14").Value = 0.04364608172438846
$ws.Range("G15").Value = 14.89209833333333
$ws.Range("H15").Value = 44.676295
$ws.Range("I15").Value = 0.07632066665966204
$ws.Range("J15").Value = 0.08031520940808551
$ws.Range("M15").Value = 1.481899
$ws.Range("N15").Value = 2.963798
$ws.Range("O15").Value = 0.00848598470259403
$ws.Range("P15").Value = 0.005673371182083786
$ws.Range("Q15").Value = 22.06858562806833
$ws.Range("R15").Value = 132.41151376841
$ws.Range("S15").Value = 0.0006476560097656703
$ws.Range("T15").Value = 0.0004556579945388569
$ws.Range("G16").Value = 14.89209833333333
$ws.Range("H16").Value = 44.676295
$ws.Range("I16").Value = 0.07632066665966204
$ws.Range("J16").Value = 0.08031520940808551
$ws.Range("M16").Value = 73.43827566666666
$ws.Range("N16").Value = 220.314827
$ws.Range("O16").Value = 0.420538838269151
$ws.Range("P16").Value = 0.4217317747321426
$ws.Range("Q16").Value = 1093.65002265844
$ws.Range("R16").Value = 9842.850203925964
$ws.Range("S16").Value = 0.0320958044929814
$ws.Range("T16").Value = 0.03387147580165557
$ws.Range("G17").Value = 29.1141605
$ws.Range("H17").Value = 58.22832099999999
$ws.Range("I17").Value = 0.1492074581338761
$ws.Range("J17").Value = 0.1046778788302885
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.546648333333333
$ws.Range("N17").Value = 7.639944999999999
$ws.Range("O17").Value = 0.01458319278139283
$ws.Range("P17").Value = 0.01462456071422719
$ws.Range("Q17").Value = 74.14352831372416
$ws.Range("R17").Value = 444.8611698823449
$ws.Range("S17").Value = 0.002175921126387914
$ws.Range("T17").Value = 0.00153086799439007
$ws.Range("G18").Value = 29.1141605
$ws.Range("H18").Value = 58.22832099999999
$ws.Range("I18").Value = 0.1492074581338761
$ws.Range("J18").Value = 0.1046778788302885
$ws.Range("O18").Value = 0.01449435301115475
$ws.Range("P18").Value = 0.01453546893349296
$ws.Range("Q18").Value = 73.69185122772717
$ws.Range("R18").Value = 442.151107366363
$ws.Range("S18").Value = 0.002162665570089493
$ws.Range("T18").Value = 0.001521542055761599
$ws.Range("G19").Value = 29.1141605
$ws.Range("H19").Value = 58.22832099999999
$ws.Range("I19").Value = 0.1492074581338761
$ws.Range("J19").Value = 0.1046778788302885
$ws.Range("M19").Value = 94.63104
$ws.Range("N19").Value = 283.89312
$ws.Range("O19").Value = 0.5418976312357076
$ws.Range("P19").Value = 0.5434348244380536
$ws.Range("Q19").Value = 2755.10328684192
$ws.Range("R19").Value = 16530.61972105152
$ws.Range("S19").Value = 0.08085516812544848
$ws.Range("T19").Value = 0.05688560470468566
$ws.Range("G20").Value = 29.1141605
$ws.Range("H20").Value = 58.22832099999999
$ws.Range("I20").Value = 0.1492074581338761
$ws.Range("J20").Value = 0.1046778788302885
$ws.Range("M20").Value = 1.481899
$ws.Range("N20").Value = 2.963798
$ws.Range("O20").Value = 0.00848598470259403
$ws.Range("P20").Value = 0.005673371182083786
$ws.Range("Q20").Value = 43.1442453307895
$ws.Range("R20").Value = 172.576981323158
$ws.Range("S20").Value = 0.001266172207237012
$ws.Range("T20").Value = 0.0005938764611574171
$ws.Range("G21").Value = 29.1141605
$ws.Range("H21").Value = 58.22832099999999
$ws.Range("I21").Value = 0.1492074581338761
$ws.Range("J21").Value = 0.1046778788302885
$ws.Range("M21").Value = 73.43827566666666
$ws.Range("N21").Value = 220.314827
$ws.Range("O21").Value = 0.420538838269151
$ws.Range("P21").Value = 0.4217317747321426
$ws.Range("Q21").Value = 2138.093744602577
$ws.Range("R21").Value = 12828.56246761546
$ws.Range("S21").Value = 0.06274753110471325
$ws.Range("T21").Value = 0.04414598761429374
$ws.Range("G22").Value = 1.32186
$ws.Range("H22").Value = 3.96558
$ws.Range("I22").Value = 0.006774413798015763
$ws.Range("J22").Value = 0.007128979431363227
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 2.546648333333333
$ws.Range("N22").Value = 7.639944999999999
$ws.Range("O22").Value = 0.01458319278139283
$ws.Range("P22").Value = 0.01462456071422719
$ws.Range("Q22").Value = 3.3663125659
$ws.Range("R22").Value = 30.2968130931
$ws.Range("S22").Value = 0.00009879258239739142
$ws.Range("T22").Value = 0.0001042581925244483
$ws.Range("G23").Value = 1.32186
$ws.Range("H23").Value = 3.96558
$ws.Range("I23").Value = 0.006774413798015763
$ws.Range("J23").Value = 0.007128979431363227
$ws.Range("O23").Value = 0.01449435301115475
$ws.Range("P23").Value = 0.01453546893349296
$ws.Range("Q23").Value = 3.34580522986
$ws.Range("R23").Value = 30.11224706874
$ws.Range("S23").Value = 0.00009819074503207803
$ws.Range("T23").Value = 0.0001036230590520905
$ws.Range("G24").Value = 1.32186
$ws.Range("H24").Value = 3.96558
$ws.Range("I24").Value = 0.006774413798015763
$ws.Range("J24").Value = 0.007128979431363227
$ws.Range("M24").Value = 94.63104
$ws.Range("N24").Value = 283.89312
$ws.Range("O24").Value = 0.5418976312357076
$ws.Range("P24").Value = 0.5434348244380536
$ws.Range("Q24").Value = 125.0889865344
$ws.Range("R24").Value = 1125.8008788096
$ws.Range("S24").Value = 0.003671038790155235
$ws.Range("T24").Value = 0.00387413568570537
$ws.Range("G25").Value = 1.32186
$ws.Range("H25").Value = 3.96558
$ws.Range("I25").Value = 0.006774413798015763
$ws.Range("J25").Value = 0.007128979431363227
$ws.Range("M25").Value = 1.481899
$ws.Range("N25").Value = 2.963798
$ws.Range("O25").Value = 0.00848598470259403
$ws.Range("P25").Value = 0.005673371182083786
$ws.Range("Q25").Value = 1.95886301214
$ws.Range("R25").Value = 11.75317807284
$ws.Range("S25").Value = 0.00005748757185900369
$ws.Range("T25").Value = 0.00004044534646356419
$ws.Range("G26").Value = 1.32186
$ws.Range("H26").Value = 3.96558
$ws.Range("I26").Value = 0.006774413798015763
$ws.Range("J26").Value = 0.007128979431363227
$ws.Range("M26").Value = 73.43827566666666
$ws.Range("N26").Value = 220.314827
$ws.Range("O26").Value = 0.420538838269151
$ws.Range("P26").Value = 0.4217317747321426
$ws.Range("Q26").Value = 97.07511907273999
$ws.Range("R26").Value = 873.6760716546599
$ws.Range("S26").Value = 0.002848904108572056
$ws.Range("T26").Value = 0.003006517147617754

Write-Output "Applied 288 cell updates"
